$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Change 1: Add a new "Make the default egg size slightly bigger." bullet
# right before the "An actual material/effect for when an egg is delivered"
# bullet (same ListParagraph / numId 25 list as its neighbour).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("An actual material/effect for when an egg is delivered", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetPara1 = $rng1.Paragraphs(1)
$xml1 = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='25'/></w:numPr></w:pPr><w:r><w:t>Make the default egg size slightly bigger.</w:t></w:r></w:p><w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='25'/></w:numPr></w:pPr><w:r><w:t>An actual material/effect for when an egg is delivered</w:t></w:r></w:p>"
$targetPara1.Range.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: Drop the stray <w:lastRenderedPageBreak/> that precedes the
# "Arenas" Heading 2 run (the first "Arenas" heading in the document).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Arenas", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetPara2 = $rng2.Paragraphs(1)
$xml2 = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>Arenas</w:t></w:r></w:p>"
$targetPara2.Range.InsertXML($xml2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: Replace the "Of course ... on top of a sombrero!?" paragraph with
# two new bold-led paragraphs describing the actual arena plan.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Of course", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetPara3 = $rng3.Paragraphs(1)
$xml3 = "<w:p xmlns:w='$wNs'><w:pPr><w:ind w:left='708'/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>YES, THIS IS THE PLAN:</w:t></w:r><w:r><w:t xml:space='preserve'> Sombrero for ground, Egg-shaped dome for boundaries</w:t></w:r></w:p><w:p xmlns:w='$wNs'><w:pPr><w:ind w:left='708'/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Arena Idea:</w:t></w:r><w:r><w:t xml:space='preserve'> two eggs next to each other =&gt; one spawns the eggs, then you use a bridge to walk to the other, to deliver it</w:t></w:r></w:p>"
$targetPara3.Range.InsertXML($xml3) | Out-Null
